$wb = $excel.ActiveWorkbook

# Rename the "veda" sheet to "Sheet1"
$vedaSheet = $wb.Worksheets.Item("veda")
$vedaSheet.Name = "Sheet1"

# Move the active tab/selection from "Sheet1" (formerly "veda") to "Sheet2",
# selecting cell H4 there.
$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet2.Activate()
$sheet2.Range("H4").Select()
